$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a third column "C" that is the sum of columns A and B for each of
# the 11 existing rows.
$ws.Range("C1").Value = 195
$ws.Range("C2").Value = 3394
$ws.Range("C3").Value = 18068
$ws.Range("C4").Value = 679
$ws.Range("C5").Value = 335
$ws.Range("C6").Value = 95585
$ws.Range("C7").Value = 7
$ws.Range("C8").Value = 7537
$ws.Range("C9").Value = 7775
$ws.Range("C10").Value = 965
$ws.Range("C11").Value = 3668

# Update the selected cell to match the new active selection in the sheet.
$ws.Range("E9").Select()
